# Apply "add data and retrospective of sprint 9" edit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Update hours data for team members (rows 5-10) ---
$ws.Range("D5").Value = 11

$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 5

$ws.Range("I7").Value = 18

$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 10

$ws.Range("D9").Value = 2
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 12

$ws.Range("C10").Value = 21

# --- Retrospective table (rows 22-29) ---
# Fill in the "BAD ESTIMATED" answers for each story.
# (Shared-string table entries get created in the order the distinct text
# values are first written, so "no" is entered before the two longer,
# unique remarks to reproduce the expected shared-string ordering.)
$ws.Range("C23").Value = "no"
$ws.Range("B23").Value = 3

$ws.Range("C24").Value = "no"
$ws.Range("C25").Value = "no"
$ws.Range("C26").Value = "no"
$ws.Range("C27").Value = "no"
$ws.Range("C28").Value = "no"
$ws.Range("C29").Value = "no"

$ws.Range("D22").Value = "required a lot of effort unexpected, in particular for the grpc connection"
$ws.Range("C22").Value = "yes, underestimated, new_func"

# --- Column E width ---
# Column E is auto-fit (bestFit) to accommodate the new, longer content added
# to the retrospective table; re-run AutoFit so the column widens accordingly.
$ws.Columns("E").AutoFit()

# --- Sheet view changes ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("C23").Select()
